$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Gaizka"
$ws.Range("C3").Value = "Detalles arreglado, creación medio hecho"
$ws.Range("D3").Value = 45765
$ws.Range("D3").NumberFormat = $ws.Range("D2").NumberFormat

$ws.Range("C3:C3").Select()
